$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 314.2
$ws.Range("I12").Value = 299
$ws.Range("J12").Value = 318
$ws.Range("K12").Value = 299
$ws.Range("L12").Value = 318
$ws.Range("M12").Value = -129
$ws.Range("N12").Value = -658

$ws.Range("H40").Value = 1724.75
$ws.Range("I40").Value = 1850
$ws.Range("K40").Value = 1850
$ws.Range("M40").Value = -1675

$ws.Range("H70").Value = 1158
$ws.Range("J70").Value = 1287.5
$ws.Range("L70").Value = 3862.5
$ws.Range("N70").Value = -4402.5

$ws.Range("H73").Value = 1158
$ws.Range("J73").Value = 1287.5
$ws.Range("L73").Value = 3862.5
$ws.Range("N73").Value = -5734.5

$ws.Range("H80").Value = 1816.3334
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1816.3334
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5449.0002
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -7445.0002

$ws.Range("H83").Value = 1816.3334
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1816.3334
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 16347.0006
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -26331.0006

$ws.Range("H132").Value = 4589.5
$ws.Range("I132").Value = 4348.45
$ws.Range("K132").Value = 13045.35
$ws.Range("M132").Value = -10515.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2717.75
$ws.Range("I45").Value = 2812.5
$ws.Range("K45").Value = 2812.5
$ws.Range("M45").Value = -2435.5

$ws.Range("H74").Value = 3266.6667
$ws.Range("I74").Value = 2400
$ws.Range("K74").Value = 2400
$ws.Range("M74").Value = -1526

$ws.Range("H77").Value = 3266.6667
$ws.Range("I77").Value = 2400
$ws.Range("K77").Value = 12000
$ws.Range("M77").Value = -7632

$ws.Range("H110").Value = 734.75
$ws.Range("I110").Value = 734.75
$ws.Range("K110").Value = 734.75
$ws.Range("M110").Value = 1310.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 87525120
$ws.Range("I7").Value = 87525120
$ws.Range("K7").Value = 87525120
$ws.Range("M7").Value = -87525007

$ws.Range("H29").Value = 9975.833000000001
$ws.Range("I29").Value = 10239
$ws.Range("K29").Value = 10239
$ws.Range("M29").Value = -9950

$ws.Range("H54").Value = 46100
$ws.Range("I54").Value = 46100
$ws.Range("K54").Value = 46100
$ws.Range("M54").Value = -45616

$ws.Range("H56").Value = 68518.336
$ws.Range("J56").Value = 75000
$ws.Range("L56").Value = 75000
$ws.Range("N56").Value = -76478

$ws.Range("H134").Value = 1956.9375
$ws.Range("I134").Value = 1960.6
$ws.Range("K134").Value = 5881.799999999999
$ws.Range("M134").Value = -3346.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 255.5
$ws.Range("I22").Value = 260
$ws.Range("K22").Value = 260
$ws.Range("M22").Value = 90

$ws.Range("H58").Value = 2892.3333
$ws.Range("I58").Value = 2869.8
$ws.Range("J58").Value = 3005
$ws.Range("K58").Value = 2869.8
$ws.Range("L58").Value = 3005
$ws.Range("M58").Value = -2666.8
$ws.Range("N58").Value = -3411

$ws.Range("H60").Value = 29999
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H122").Value = 1500
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 2892.3333
$ws.Range("I136").Value = 2869.8
$ws.Range("J136").Value = 3005
$ws.Range("K136").Value = 8609.400000000001
$ws.Range("L136").Value = 9015
$ws.Range("M136").Value = -6059.400000000001
$ws.Range("N136").Value = -14115

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 895.2222
$ws.Range("J38").Value = 82.59999999999999
$ws.Range("L38").Value = 247.8
$ws.Range("N38").Value = -941.8

$ws.Range("H68").Value = 461.23077
$ws.Range("I68").Value = 503.7143
$ws.Range("J68").Value = 411.66666
$ws.Range("K68").Value = 1511.1429
$ws.Range("L68").Value = 1234.99998
$ws.Range("M68").Value = -700.1428999999998
$ws.Range("N68").Value = -2856.99998

$ws.Range("H71").Value = 461.23077
$ws.Range("I71").Value = 503.7143
$ws.Range("J71").Value = 411.66666
$ws.Range("K71").Value = 4533.428699999999
$ws.Range("L71").Value = 3704.99994
$ws.Range("M71").Value = -477.4286999999995
$ws.Range("N71").Value = -11816.99994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 23333.334
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -9404
$ws.Range("N44").Value = -31192

$ws.Range("H47").Value = 29950

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H123").Value = 98999
$ws.Range("J123").Value = 98999
$ws.Range("L123").Value = 98999
$ws.Range("N123").Value = -103899

$ws.Range("H127").Value = 99999
$ws.Range("J127").Value = 99999
$ws.Range("L127").Value = 99999
$ws.Range("N127").Value = -109919

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888

$ws.Range("H16").Value = 2156.75
$ws.Range("I16").Value = 2209.7334
$ws.Range("J16").Value = 1997.8
$ws.Range("K16").Value = 2209.7334
$ws.Range("L16").Value = 1997.8
$ws.Range("M16").Value = -2039.7334
$ws.Range("N16").Value = -2337.8

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H22").Value = 1500
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2090

$ws.Range("H27").Value = 1500
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1714

$ws.Range("H123").Value = 27500
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws.Range("H127").Value = 20000
$ws.Range("J127").Value = 20000
$ws.Range("L127").Value = 20000
$ws.Range("N127").Value = -29920

$ws.Range("H136").Value = 11011.2
$ws.Range("I136").Value = 2525.75
$ws.Range("K136").Value = 7577.25
$ws.Range("M136").Value = -5027.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -36108

$ws.Range("H62").Value = 2204.6
$ws.Range("J62").Value = 2341
$ws.Range("L62").Value = 2341
$ws.Range("N62").Value = -3589

$ws.Range("H65").Value = 2204.6
$ws.Range("J65").Value = 2341
$ws.Range("L65").Value = 11705
$ws.Range("N65").Value = -17945

$ws.Range("H81").Value = 2000.5
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 2
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = 1059
$ws.Range("N81").Value = -10122

$ws.Range("H84").Value = 2000.5
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 10
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = 5294
$ws.Range("N84").Value = -50608

$ws.Range("H132").Value = 1998.6
$ws.Range("I132").Value = 1331.3334
$ws.Range("K132").Value = 3994.0002
$ws.Range("M132").Value = -1464.0002

$ws.Range("H136").Value = 2882
$ws.Range("I136").Value = 2646.6667
$ws.Range("K136").Value = 7940.000100000001
$ws.Range("M136").Value = -5390.000100000001
